$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures to the leve-profit
# tables on several job sheets, as produced by the scheduled price-refresh runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3676.8533
$ws.Range("I64").Value = 3441.4883
$ws.Range("J64").Value = 3993.125
$ws.Range("K64").Value = 3441.4883
$ws.Range("L64").Value = 3993.125
$ws.Range("M64").Value = -3193.4883
$ws.Range("N64").Value = -4489.125
$ws.Range("H67").Value = 3676.8533
$ws.Range("I67").Value = 3441.4883
$ws.Range("J67").Value = 3993.125
$ws.Range("K67").Value = 3441.4883
$ws.Range("L67").Value = 3993.125
$ws.Range("M67").Value = -2583.4883
$ws.Range("N67").Value = -5709.125
$ws.Range("H74").Value = 3515.2778
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 3471.2122
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 3471.2122
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -5343.2122
$ws.Range("H77").Value = 3515.2778
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 3471.2122
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 17356.061
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -26716.061
$ws.Range("H132").Value = 2899.1177
$ws.Range("I132").Value = 2368.7437
$ws.Range("J132").Value = 4622.8335
$ws.Range("K132").Value = 7106.2311
$ws.Range("L132").Value = 13868.5005
$ws.Range("M132").Value = -4576.2311
$ws.Range("N132").Value = -18928.5005
$ws.Range("H137").Value = 12989.721
$ws.Range("I137").Value = 769.9048
$ws.Range("J137").Value = 46461.39
$ws.Range("K137").Value = 2309.7144
$ws.Range("L137").Value = 139384.17
$ws.Range("M137").Value = 240.2856000000002
$ws.Range("N137").Value = -144484.17
$ws.Range("H138").Value = 3244.8298
$ws.Range("I138").Value = 1293.6
$ws.Range("K138").Value = 3880.8
$ws.Range("M138").Value = 1259.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10740.012
$ws.Range("I32").Value = 4973.3623
$ws.Range("J32").Value = 37266.6
$ws.Range("K32").Value = 4973.3623
$ws.Range("L32").Value = 37266.6
$ws.Range("M32").Value = -4686.3623
$ws.Range("N32").Value = -37840.6
$ws.Range("H76").Value = 20288
$ws.Range("J76").Value = 20288
$ws.Range("L76").Value = 20288
$ws.Range("N76").Value = -20964
$ws.Range("H79").Value = 20288
$ws.Range("J79").Value = 20288
$ws.Range("L79").Value = 20288
$ws.Range("N79").Value = -22628
$ws.Range("H132").Value = 1693.6428
$ws.Range("I132").Value = 806
$ws.Range("J132").Value = 2717.8462
$ws.Range("K132").Value = 2418
$ws.Range("L132").Value = 8153.5386
$ws.Range("M132").Value = 112
$ws.Range("N132").Value = -13213.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1421.6562
$ws.Range("I86").Value = 1513.7333
$ws.Range("J86").Value = 1340.4117
$ws.Range("K86").Value = 1513.7333
$ws.Range("L86").Value = 1340.4117
$ws.Range("M86").Value = -390.7333000000001
$ws.Range("N86").Value = -3586.4117
$ws.Range("H89").Value = 1421.6562
$ws.Range("I89").Value = 1513.7333
$ws.Range("J89").Value = 1340.4117
$ws.Range("K89").Value = 7568.6665
$ws.Range("L89").Value = 6702.058500000001
$ws.Range("M89").Value = -1952.6665
$ws.Range("N89").Value = -17934.0585
$ws.Range("H107").Value = 1394
$ws.Range("I107").Value = 1314.4286
$ws.Range("J107").Value = 1505.4
$ws.Range("K107").Value = 1314.4286
$ws.Range("L107").Value = 1505.4
$ws.Range("M107").Value = 605.5714
$ws.Range("N107").Value = -5345.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1992.8959
$ws.Range("I31").Value = 1465.2565
$ws.Range("J31").Value = 4279.3335
$ws.Range("K31").Value = 1465.2565
$ws.Range("L31").Value = 4279.3335
$ws.Range("M31").Value = -1170.2565
$ws.Range("N31").Value = -4869.3335
$ws.Range("H34").Value = 1992.8959
$ws.Range("I34").Value = 1465.2565
$ws.Range("J34").Value = 4279.3335
$ws.Range("K34").Value = 1465.2565
$ws.Range("L34").Value = 4279.3335
$ws.Range("M34").Value = -1263.2565
$ws.Range("N34").Value = -4683.3335
$ws.Range("H62").Value = 2418.4546
$ws.Range("I62").Value = 2371.8572
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2371.8572
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1747.8572
$ws.Range("N62").Value = -3748
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 35271
$ws.Range("J64").Value = 35271
$ws.Range("L64").Value = 35271
$ws.Range("N64").Value = -35767
$ws.Range("H65").Value = 2418.4546
$ws.Range("I65").Value = 2371.8572
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 11859.286
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -8739.286
$ws.Range("N65").Value = -18740
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 35271
$ws.Range("J67").Value = 35271
$ws.Range("L67").Value = 35271
$ws.Range("N67").Value = -36987
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 21000
$ws.Range("J70").Value = 21000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21630
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 21000
$ws.Range("J73").Value = 21000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -23184
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H82").Value = 16998
$ws.Range("I82").Value = 1994
$ws.Range("K82").Value = 1994
$ws.Range("M82").Value = -1633
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H85").Value = 16998
$ws.Range("I85").Value = 1994
$ws.Range("K85").Value = 1994
$ws.Range("M85").Value = -746
$ws.Range("H88").Value = 20343
$ws.Range("J88").Value = 20343
$ws.Range("L88").Value = 20343
$ws.Range("N88").Value = -21155
$ws.Range("H91").Value = 20343
$ws.Range("J91").Value = 20343
$ws.Range("L91").Value = 20343
$ws.Range("N91").Value = -23151
$ws.Range("H134").Value = 2382.6167
$ws.Range("I134").Value = 1393.5625
$ws.Range("J134").Value = 3512.9644
$ws.Range("K134").Value = 4180.6875
$ws.Range("L134").Value = 10538.8932
$ws.Range("M134").Value = -1645.6875
$ws.Range("N134").Value = -15608.8932
$ws.Range("H141").Value = 39646.668
$ws.Range("J141").Value = 39646.668
$ws.Range("L141").Value = 39646.668
$ws.Range("N141").Value = -50006.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2562
$ws.Range("I69").Value = 1966.3334
$ws.Range("J69").Value = 2859.8333
$ws.Range("K69").Value = 5899.0002
$ws.Range("L69").Value = 8579.499899999999
$ws.Range("M69").Value = -5088.0002
$ws.Range("N69").Value = -10201.4999
$ws.Range("H72").Value = 2562
$ws.Range("I72").Value = 1966.3334
$ws.Range("J72").Value = 2859.8333
$ws.Range("K72").Value = 17697.0006
$ws.Range("L72").Value = 25738.4997
$ws.Range("M72").Value = -13641.0006
$ws.Range("N72").Value = -33850.4997
$ws.Range("H74").Value = 9155.714
$ws.Range("I74").Value = 390
$ws.Range("J74").Value = 10616.667
$ws.Range("K74").Value = 1170
$ws.Range("L74").Value = 31850.001
$ws.Range("M74").Value = -109
$ws.Range("N74").Value = -33972.001
$ws.Range("H77").Value = 9155.714
$ws.Range("I77").Value = 390
$ws.Range("J77").Value = 10616.667
$ws.Range("K77").Value = 3510
$ws.Range("L77").Value = 95550.003
$ws.Range("M77").Value = 1794
$ws.Range("N77").Value = -106158.003
$ws.Range("H80").Value = 5000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16872
$ws.Range("H83").Value = 5000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54360
$ws.Range("H113").Value = 442.78262
$ws.Range("I113").Value = 408
$ws.Range("J113").Value = 469.53845
$ws.Range("K113").Value = 1224
$ws.Range("L113").Value = 1408.61535
$ws.Range("M113").Value = 946
$ws.Range("N113").Value = -5748.61535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2831.9143
$ws.Range("I80").Value = 2661.8572
$ws.Range("J80").Value = 3087
$ws.Range("K80").Value = 2661.8572
$ws.Range("L80").Value = 3087
$ws.Range("M80").Value = -1663.8572
$ws.Range("N80").Value = -5083
$ws.Range("H83").Value = 2831.9143
$ws.Range("I83").Value = 2661.8572
$ws.Range("J83").Value = 3087
$ws.Range("K83").Value = 13309.286
$ws.Range("L83").Value = 15435
$ws.Range("M83").Value = -8317.286
$ws.Range("N83").Value = -25419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 858.875
$ws.Range("I46").Value = 397.75
$ws.Range("J46").Value = 1320
$ws.Range("K46").Value = 397.75
$ws.Range("L46").Value = 1320
$ws.Range("M46").Value = -209.75
$ws.Range("N46").Value = -1696
$ws.Range("H55").Value = 198.73077
$ws.Range("I55").Value = 86.5
$ws.Range("J55").Value = 294.92856
$ws.Range("K55").Value = 86.5
$ws.Range("L55").Value = 294.92856
$ws.Range("M55").Value = 86.5
$ws.Range("N55").Value = -640.9285600000001
$ws.Range("H132").Value = 2430.05
$ws.Range("I132").Value = 1989.9131
$ws.Range("J132").Value = 3409.7097
$ws.Range("K132").Value = 5969.7393
$ws.Range("L132").Value = 10229.1291
$ws.Range("M132").Value = -3439.7393
$ws.Range("N132").Value = -15289.1291
